$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.00469825789230973
$ws.Range("J2").Value = 0.00469825789230973
$ws.Range("M2").Value = 71.05094633333333
$ws.Range("N2").Value = 213.152839
$ws.Range("O2").Value = 0.8240565632932695
$ws.Range("P2").Value = 0.8240565632932696
$ws.Range("Q2").Value = 6.638171497679444
$ws.Range("R2").Value = 59.743543479115
$ws.Range("S2").Value = 0.003871630252202236
$ws.Range("T2").Value = 0.003871630252202237
$ws.Range("I3").Value = 0.00469825789230973
$ws.Range("J3").Value = 0.00469825789230973
$ws.Range("O3").Value = 0.1323102827659759
$ws.Range("P3").Value = 0.132310282765976
$ws.Range("S3").Value = 0.0006216278302389785
$ws.Range("T3").Value = 0.0006216278302389786
$ws.Range("I4").Value = 0.00469825789230973
$ws.Range("J4").Value = 0.00469825789230973
$ws.Range("M4").Value = 3.762092333333333
$ws.Range("N4").Value = 11.286277
$ws.Range("O4").Value = 0.04363315394075456
$ws.Range("P4").Value = 0.04363315394075455
$ws.Range("Q4").Value = 0.3514860165494444
$ws.Range("R4").Value = 3.163374148945
$ws.Range("S4").Value = 0.0002049998098685155
$ws.Range("T4").Value = 0.0002049998098685155
$ws.Range("I5").Value = 0.7185612021237531
$ws.Range("J5").Value = 0.7185612021237531
$ws.Range("M5").Value = 71.05094633333333
$ws.Range("N5").Value = 213.152839
$ws.Range("O5").Value = 0.8240565632932695
$ws.Range("P5").Value = 0.8240565632932696
$ws.Range("Q5").Value = 1015.255569321506
$ws.Range("R5").Value = 9137.300123893556
$ws.Range("S5").Value = 0.5921350747379803
$ws.Range("T5").Value = 0.5921350747379804
$ws.Range("I6").Value = 0.7185612021237531
$ws.Range("J6").Value = 0.7185612021237531
$ws.Range("O6").Value = 0.1323102827659759
$ws.Range("P6").Value = 0.132310282765976
$ws.Range("S6").Value = 0.09507303583765335
$ws.Range("T6").Value = 0.09507303583765336
$ws.Range("I7").Value = 0.7185612021237531
$ws.Range("J7").Value = 0.7185612021237531
$ws.Range("M7").Value = 3.762092333333333
$ws.Range("N7").Value = 11.286277
$ws.Range("O7").Value = 0.04363315394075456
$ws.Range("P7").Value = 0.04363315394075455
$ws.Range("Q7").Value = 53.75699256417233
$ws.Range("R7").Value = 483.8129330775509
$ws.Range("S7").Value = 0.03135309154811937
$ws.Range("T7").Value = 0.03135309154811936
$ws.Range("G8").Value = 5.503190333333333
$ws.Range("H8").Value = 16.509571
$ws.Range("I8").Value = 0.2767405399839373
$ws.Range("J8").Value = 0.2767405399839373
$ws.Range("M8").Value = 71.05094633333333
$ws.Range("N8").Value = 213.152839
$ws.Range("O8").Value = 0.8240565632932695
$ws.Range("P8").Value = 0.8240565632932696
$ws.Range("Q8").Value = 391.0068810357854
$ws.Range("R8").Value = 3519.061929322069
$ws.Range("S8").Value = 0.228049858303087
$ws.Range("T8").Value = 0.228049858303087
$ws.Range("G9").Value = 5.503190333333333
$ws.Range("H9").Value = 16.509571
$ws.Range("I9").Value = 0.2767405399839373
$ws.Range("J9").Value = 0.2767405399839373
$ws.Range("O9").Value = 0.1323102827659759
$ws.Range("P9").Value = 0.132310282765976
$ws.Range("Q9").Value = 62.77995139864645
$ws.Range("R9").Value = 565.0195625878181
$ws.Range("S9").Value = 0.03661561909808361
$ws.Range("T9").Value = 0.03661561909808361
$ws.Range("G10").Value = 5.503190333333333
$ws.Range("H10").Value = 16.509571
$ws.Range("I10").Value = 0.2767405399839373
$ws.Range("J10").Value = 0.2767405399839373
$ws.Range("M10").Value = 3.762092333333333
$ws.Range("N10").Value = 11.286277
$ws.Range("O10").Value = 0.04363315394075456
$ws.Range("P10").Value = 0.04363315394075455
$ws.Range("Q10").Value = 20.70351016190744
$ws.Range("R10").Value = 186.331591457167
$ws.Range("S10").Value = 0.01207506258276668
$ws.Range("T10").Value = 0.01207506258276667
